# Fix the two mis-sorted HeapSort timing samples on the "Data" sheet
# (commit: "Fix sorting and generate viable xlsx and charts").
#
#   D2 (Avg_Time_ms for the 5000-element run):  0.97829543 -> 0.84116882
#   D3 (Avg_Time_ms for the 10000-element run): 1.7768945  -> 1.8166705
#
# The sheet's scatter chart plots Data!$D$2:$D$8 / Data!$E$2:$E$8 directly,
# so correcting these two source cells is the authoritative fix - the
# chart re-reads its data from the worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$ws.Range("D2").Value = 0.84116882
$ws.Range("D3").Value = 1.8166705

# Recalculate so every dependent (formulas, chart caches) picks up the
# corrected inputs.
$excel.CalculateFullRebuild()
